# Refresh the crypto price/volume figures (columns D and E) for rows 2-51
# to match the latest scrape, as produced by the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.977.02"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "'1.883.03"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("D4").Value = "'0.9989"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'243.37"
$ws.Range("E5").Value = "  -3.33%  "
$ws.Range("D6").Value = "'0.9986"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "'0.4939"
$ws.Range("E7").Value = "  -3.29%  "
$ws.Range("D8").Value = "'0.2945"
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").Value = "'0.06647"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("D10").Value = "'1.880.64"
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("D11").Value = "'16.77"
$ws.Range("E11").Value = "  -2.98%  "
$ws.Range("D12").Value = "'0.07203"
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("D13").Value = "'0.6690"
$ws.Range("E13").Value = "  -4.32%  "
$ws.Range("D14").Value = "'86.43"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").Value = "'29.948.13"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "'0.000007850"
$ws.Range("E17").Value = "  -4.10%  "
$ws.Range("D18").Value = "'0.9985"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("D20").Value = "'2.119.26"
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("D21").Value = "'0.9990"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "'4.787"
$ws.Range("D23").Value = "'5.879"
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("D24").Value = "'9.109"
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("D25").Value = "'150.50"
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("D26").Value = "'142.62"
$ws.Range("E26").Value = "  +5.46%  "
$ws.Range("D27").Value = "'17.09"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").Value = "'1.925"
$ws.Range("E28").Value = "  -3.59%  "
$ws.Range("D29").Value = "'1.389"
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("D30").Value = "'4.220"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").Value = "'0.08792"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").Value = "'4.018"
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").Value = "'0.05058"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "'0.7149"
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("D35").Value = "'1.117"
$ws.Range("E35").Value = "  -2.16%  "
$ws.Range("D36").Value = "'2.666"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("E37").Value = "  +5.83%  "
$ws.Range("D38").Value = "'2.695"
$ws.Range("E38").Value = "  -4.25%  "
$ws.Range("D39").Value = "'2.174"
$ws.Range("E39").Value = "  -4.02%  "
$ws.Range("D40").Value = "'0.9321"
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("D41").Value = "'5.767"
$ws.Range("E41").Value = "  -6.27%  "
$ws.Range("D42").Value = "'0.4238"
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("D43").Value = "'0.9986"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "'103.12"
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("D45").Value = "'7.416"
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("D46").Value = "'0.1274"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").Value = "'0.05672"
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("D48").Value = "'32.60"
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("D49").Value = "'8.324"
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("D50").Value = "'0.3776"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "'56.17"
$ws.Range("E51").Value = "  -1.33%  "
